# Update the "2019" sheet (sheet1.xml / rId1) for the Day 5 entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# Row 9 was the placeholder "Day 5: TITLE" row with no recorded times.
# It now holds the real Day 5 puzzle name and its timing data.
$ws.Range("B9").Value = "Day 5: Sunny with a Chance of Asteroids"
$ws.Range("C9").Value = 0.019201388888888889
$ws.Range("E9").Value = 0.029027777777777777
$ws.Range("F9").Value = 0.012395833333333335
$ws.Range("H9").Value = "3rd"

# Update the active cell/selection to match the saved view state.
$ws.Range("K13").Select()

$wb.Save()
